# Apply crypto price/volume updates (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string need to be
# pre-formatted as Text so COM Range.Value assignment does not coerce
# them into a Double (the source workbook stores every Price/Volume
# cell as a literal string, even when it looks like a number).
$textCells = @("D5", "D9", "D10", "D15", "D16", "D18", "D23", "D25", "D27", "D35", "D36", "D37", "D38", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.955.59"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.672.15"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "214.98"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.0620"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "20.16"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "1.908.06"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "1.678.15"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "0.525"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "65.59"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "26.951.97"
$ws.Range("D18").Value = "234.54"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").Value = "0.0₃0732"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "9.15"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").Value = "145.75"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "15.96"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "1.470.43"
$ws.Range("E33").Value = "  -5.40%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").Value = "1.65"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "0.578"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "0.896"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +8.17%  "
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D43").Value = "2.30"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").Value = "66.68"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "1.813.72"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").Value = "0.779"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "90.51"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "1.54"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "7.70"

# Restore the default style on the cells we temporarily reformatted,
# so only the cell contents change (matches the source diff, which
# shows no styling changes).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
